# #temp fixed antiforgery exception
# Update the sample/test GUID used for the "WResponses" seed row so that it
# no longer collides with the other fixture id. Also move the sheet's
# active selection as it was left after making the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WResponses")

# B3 held the placeholder response id "d6485221-65g8-42a0-9873-00622ec70e79"
# (sharedStrings index 22) - replace it with a fresh guid-shaped value.
$ws.Range("B3").Value = "d6485221-63e8-42a0-9873-00622ec70e79"

# Leave the selection where the editor ended up (B5) instead of the old B12.
$ws.Range("B5").Select() | Out-Null
